$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 224.7833046666667
$ws.Range("H2").Value = 674.349914
$ws.Range("I2").Value = 0.3882379172278888
$ws.Range("J2").Value = 0.3882379172278889
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 32751.09487856755
$ws.Range("R2").Value = 294759.8539071079
$ws.Range("S2").Value = 0.1112665703005874
$ws.Range("T2").Value = 0.1112665703005874

$ws.Range("G3").Value = 224.7833046666667
$ws.Range("H3").Value = 674.349914
$ws.Range("I3").Value = 0.3882379172278888
$ws.Range("J3").Value = 0.3882379172278889
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 37943.35499216409
$ws.Range("R3").Value = 341490.1949294768
$ws.Range("S3").Value = 0.1289064378253366
$ws.Range("T3").Value = 0.1289064378253367

$ws.Range("G4").Value = 224.7833046666667
$ws.Range("H4").Value = 674.349914
$ws.Range("I4").Value = 0.3882379172278888
$ws.Range("J4").Value = 0.3882379172278889
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 28800.62046020579
$ws.Range("R4").Value = 259205.5841418521
$ws.Range("S4").Value = 0.09784546968636131
$ws.Range("T4").Value = 0.09784546968636132

$ws.Range("G5").Value = 224.7833046666667
$ws.Range("H5").Value = 674.349914
$ws.Range("I5").Value = 0.3882379172278888
$ws.Range("J5").Value = 0.3882379172278889
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 14781.99265606573
$ws.Range("R5").Value = 133037.9339045916
$ws.Range("S5").Value = 0.05021943941560349
$ws.Range("T5").Value = 0.0502194394156035

$ws.Range("I6").Value = 0.4251955538547045
$ws.Range("J6").Value = 0.4251955538547046
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 35868.77867487222
$ws.Range("R6").Value = 322819.00807385
$ws.Range("S6").Value = 0.1218583988969359
$ws.Range("T6").Value = 0.121858398896936

$ws.Range("I7").Value = 0.4251955538547045
$ws.Range("J7").Value = 0.4251955538547046
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.1411774630822787
$ws.Range("T7").Value = 0.1411774630822787

$ws.Range("I8").Value = 0.4251955538547045
$ws.Range("J8").Value = 0.4251955538547046
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 31542.24568114044
$ws.Range("R8").Value = 283880.211130264
$ws.Range("S8").Value = 0.1071596998369574
$ws.Range("T8").Value = 0.1071596998369575

$ws.Range("I9").Value = 0.4251955538547045
$ws.Range("J9").Value = 0.4251955538547046
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 16189.13886451416
$ws.Range("R9").Value = 145702.2497806274
$ws.Range("S9").Value = 0.05499999203853245
$ws.Range("T9").Value = 0.05499999203853246

$ws.Range("G10").Value = 107.695137
$ws.Range("H10").Value = 323.085411
$ws.Range("I10").Value = 0.186007300437435
$ws.Range("J10").Value = 0.186007300437435
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 15691.26165787869
$ws.Range("R10").Value = 141221.3549209082
$ws.Range("S10").Value = 0.05330853441189239
$ws.Range("T10").Value = 0.0533085344118924

$ws.Range("G11").Value = 107.695137
$ws.Range("H11").Value = 323.085411
$ws.Range("I11").Value = 0.186007300437435
$ws.Range("J11").Value = 0.186007300437435
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 18178.90710424594
$ws.Range("R11").Value = 163610.1639382134
$ws.Range("S11").Value = 0.06175990918172614
$ws.Range("T11").Value = 0.06175990918172616

$ws.Range("G12").Value = 107.695137
$ws.Range("H12").Value = 323.085411
$ws.Range("I12").Value = 0.186007300437435
$ws.Range("J12").Value = 0.186007300437435
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 13798.56378010986
$ws.Range("R12").Value = 124187.0740209887
$ws.Range("S12").Value = 0.04687839819032895
$ws.Range("T12").Value = 0.04687839819032896

$ws.Range("G13").Value = 107.695137
$ws.Range("H13").Value = 323.085411
$ws.Range("I13").Value = 0.186007300437435
$ws.Range("J13").Value = 0.186007300437435
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 7082.148412172823
$ws.Range("R13").Value = 63739.33570955541
$ws.Range("S13").Value = 0.0240604586534875
$ws.Range("T13").Value = 0.02406045865348751

$ws.Range("G14").Value = 0.323784
$ws.Range("H14").Value = 0.971352
$ws.Range("I14").Value = 0.0005592284799715185
$ws.Range("J14").Value = 0.0005592284799715186
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 47.175569911152
$ws.Range("R14").Value = 424.580129200368
$ws.Range("S14").Value = 0.0001602714011684684
$ws.Range("T14").Value = 0.0001602714011684684

$ws.Range("G15").Value = 0.323784
$ws.Range("H15").Value = 0.971352
$ws.Range("I15").Value = 0.0005592284799715185
$ws.Range("J15").Value = 0.0005592284799715186
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 54.65464292822401
$ws.Range("R15").Value = 491.8917863540161
$ws.Range("S15").Value = 0.0001856803472425688
$ws.Range("T15").Value = 0.0001856803472425688

$ws.Range("G16").Value = 0.323784
$ws.Range("H16").Value = 0.971352
$ws.Range("I16").Value = 0.0005592284799715185
$ws.Range("J16").Value = 0.0005592284799715186
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 41.48519886259201
$ws.Range("R16").Value = 373.366789763328
$ws.Range("S16").Value = 0.0001409392819627266
$ws.Range("T16").Value = 0.0001409392819627266

$ws.Range("G17").Value = 0.323784
$ws.Range("H17").Value = 0.971352
$ws.Range("I17").Value = 0.0005592284799715185
$ws.Range("J17").Value = 0.0005592284799715186
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 21.292385202936
$ws.Range("R17").Value = 191.631466826424
$ws.Range("S17").Value = 0.0000723374495977548
$ws.Range("T17").Value = 0.00007233744959775481
